# SA Done.xlsx - update SS courses: renumber 1168 -> 1178 on sheet "1168",
# renumber 1083 -> 1084 on sheet "1171", fix a few LAB->LEC component values,
# clean up a couple of border artifacts, and refresh the view state
# (active sheet / selection) to match the saved workbook.

$wb = $excel.ActiveWorkbook

$ws1168 = $wb.Worksheets.Item("1168")
$ws1171 = $wb.Worksheets.Item("1171")

# ---------------------------------------------------------------------
# Sheet "1168": rows 3-49, column A holds the catalog number, 1168 -> 1178
# ---------------------------------------------------------------------
for ($r = 3; $r -le 49; $r++) {
    $ws1168.Cells.Item($r, 1).Value = 1178
}

# Row 23, column I (Component) was mis-entered as LAB, should be LEC
$ws1168.Range("I23").Value = "LEC"

# Rows 25/26 sit at the boundary between the two course blocks; the
# border between them should no longer be drawn (the blocks are now a
# single continuous table), so drop A25's bottom border and A26's top
# border.
$ws1168.Range("A25").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws1168.Range("A26").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none

# ---------------------------------------------------------------------
# Sheet "1171": rows 3-9, column F holds the catalog number, 1083 -> 1084
# and column I (Component) corrected from LAB to LEC
# ---------------------------------------------------------------------
for ($r = 3; $r -le 9; $r++) {
    $ws1171.Cells.Item($r, 6).Value = 1084
    $ws1171.Range("I" + $r).Value = "LEC"
}

# Row 3's "I" cell also loses its stray top border (it previously
# inherited the header-row border along with the rest of row 3).
$ws1171.Range("I3").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none

# ---------------------------------------------------------------------
# View state: sheet "1171" becomes the active/visible tab, with a new
# selection; sheet "1168" is scrolled down and also gets a new selection.
# ---------------------------------------------------------------------
$null = $ws1168.Activate()
$null = $ws1168.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws1168.Range("J30").Select()

$null = $ws1171.Activate()
$null = $ws1171.Range("L7").Select()
